# Update results of test case 001 -> 004
# Appends new running_logs rows (22-31) to the logs sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$preprocess1 = 'remove multiple spaces, convert to lower, trim "space" and ",", convert unicode to ascii'
$model1 = '2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 4000'

$preprocess2 = 'remove multiple spaces, convert to lower, convert unicode to ascii, trim "space" and ","'
$model2 = '2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 6000'

$features = '11 features: #ascii/(#ascii+#digit+#punctuation), #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, #max_digit_skip_0_1, first_character_ascii, first_character_digit, #(, #+, #/'
$model = 'Neuron Network'
$filter = '0 filters: '

$rows = @(
    @{ Row=22; Time='20160405_164907'; RunningTime=2255.573; Preprocess=$preprocess1; ModelDetails=$model1; Test=0.995333333333333; Val=0.99009900990099; J=0.364583333333333 },
    @{ Row=23; Time='20160405_172643'; RunningTime=2320.508; Preprocess=$preprocess1; ModelDetails=$model1; Test=0.994666666666667; Val=0.99009900990099; J=0.34375 },
    @{ Row=24; Time='20160405_180523'; RunningTime=2315.106; Preprocess=$preprocess1; ModelDetails=$model1; Test=0.992; Val=0.99009900990099; J=0.354166666666667 },
    @{ Row=25; Time='20160405_184358'; RunningTime=2282.16; Preprocess=$preprocess1; ModelDetails=$model1; Test=0.994; Val=0.99009900990099; J=0.3125 },
    @{ Row=26; Time='20160405_192200'; RunningTime=2303.633; Preprocess=$preprocess1; ModelDetails=$model1; Test=0.995333333333333; Val=0.99009900990099; J=0.322916666666667 },
    @{ Row=27; Time='20160406_081500'; RunningTime=3489.767; Preprocess=$preprocess2; ModelDetails=$model2; Test=0.993333333333333; Val=0.99009900990099; J=0.302083333333333 },
    @{ Row=28; Time='20160406_091310'; RunningTime=3582.836; Preprocess=$preprocess2; ModelDetails=$model2; Test=0.994; Val=0.993399339933993; J=0.278350515463918 },
    @{ Row=29; Time='20160406_101253'; RunningTime=5528.333; Preprocess=$preprocess2; ModelDetails=$model2; Test=0.994666666666667; Val=0.99009900990099; J=0.333333333333333 },
    @{ Row=30; Time='20160406_114501'; RunningTime=5596.762; Preprocess=$preprocess2; ModelDetails=$model2; Test=0.989333333333333; Val=0.99009900990099; J=0.3125 },
    @{ Row=31; Time='20160406_131818'; RunningTime=7310.106; Preprocess=$preprocess2; ModelDetails=$model2; Test=0.993333333333333; Val=0.99009900990099; J=0.270833333333333 }
)

foreach ($r in $rows) {
    $rowIndex = $r.Row
    $ws.Cells.Item($rowIndex, 1).Value = $r.Time
    $ws.Cells.Item($rowIndex, 2).Value = $r.RunningTime
    $ws.Cells.Item($rowIndex, 3).Value = $r.Preprocess
    $ws.Cells.Item($rowIndex, 4).Value = $features
    $ws.Cells.Item($rowIndex, 5).Value = $model
    $ws.Cells.Item($rowIndex, 6).Value = $r.ModelDetails
    $ws.Cells.Item($rowIndex, 7).Value = $r.Test
    $ws.Cells.Item($rowIndex, 8).Value = $r.Val
    $ws.Cells.Item($rowIndex, 9).Value = $filter
    $ws.Cells.Item($rowIndex, 10).Value = $r.J
}
